$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue {
    param($Worksheet, [string]$Address, [string]$Text)
    $cell = $Worksheet.Range($Address)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = $originalStyle
}

# Row 2
$ws.Range('D2').Value = '47.415.84'
$ws.Range('E2').Value = '  +2.91%  '

# Row 3
$ws.Range('D3').Value = '2.512.94'
$ws.Range('E3').Value = '  +2.45%  '

# Row 4
$ws.Range('E4').Value = '  +0.10%  '

# Row 5
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextCellValue $ws 'D5' '324.73'
$ws.Range('E5').Value = '  +1.25%  '

# Row 6
$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCellValue $ws 'D6' '110.07'
$ws.Range('E6').Value = '  +5.19%  '

# Row 7
$ws.Range('E7').Value = '  +1.47%  '

# Row 8
Set-TextCellValue $ws 'D8' '0.999'
$ws.Range('E8').Value = '  -0.05%  '

# Row 9
Set-TextCellValue $ws 'D9' '0.541'
$ws.Range('E9').Value = '  +1.12%  '

# Row 10
Set-TextCellValue $ws 'D10' '39.35'
$ws.Range('E10').Value = '  +9.47%  '

# Row 11
$ws.Range('E11').Value = '  +1.55%  '

# Row 13
Set-TextCellValue $ws 'D13' '18.65'
$ws.Range('E13').Value = '  +2.22%  '

# Row 14
$ws.Range('E14').Value = '  +2.51%  '

# Row 15
$ws.Range('D15').Value = '2.905.13'
$ws.Range('E15').Value = '  +2.53%  '

# Row 16
$ws.Range('D16').Value = '2.507.76'
$ws.Range('E16').Value = '  +3.03%  '

# Row 17
Set-TextCellValue $ws 'D17' '0.865'
$ws.Range('E17').Value = '  +2.95%  '

# Row 18
$ws.Range('D18').Value = '47.373.41'
$ws.Range('E18').Value = '  +3.15%  '

# Row 19
Set-TextCellValue $ws 'D19' '12.93'
$ws.Range('E19').Value = '  +2.71%  '

# Row 20
Set-TextCellValue $ws 'D20' '6.73'
$ws.Range('E20').Value = '  +4.82%  '

# Row 21
$ws.Range('D21').Value = '0.0₃0946'
$ws.Range('E21').Value = '  +1.17%  '

# Row 22
$ws.Range('E22').Value = '  +12.76%  '

# Row 23
Set-TextCellValue $ws 'D23' '70.82'
$ws.Range('E23').Value = '  -0.86%  '

# Row 24
Set-TextCellValue $ws 'D24' '249.83'
$ws.Range('E24').Value = '  +1.23%  '

# Row 25
Set-TextCellValue $ws 'D25' '2.60'
$ws.Range('E25').Value = '  +4.05%  '

# Row 26
Set-TextCellValue $ws 'D26' '26.24'
$ws.Range('E26').Value = '  +0.82%  '

# Row 27
Set-TextCellValue $ws 'D27' '1.00'
$ws.Range('E27').Value = '  -0.03%  '

# Row 28
$ws.Range('E28').Value = '  +4.54%  '

# Row 29
Set-TextCellValue $ws 'D29' '10.06'
$ws.Range('E29').Value = '  +3.55%  '

# Row 30
Set-TextCellValue $ws 'D30' '35.66'
$ws.Range('E30').Value = '  +4.80%  '

# Row 31
Set-TextCellValue $ws 'D31' '0.138'
$ws.Range('E31').Value = '  +6.27%  '

# Row 32
Set-TextCellValue $ws 'D32' '50.25'
$ws.Range('E32').Value = '  +1.59%  '

# Row 33
Set-TextCellValue $ws 'D33' '19.99'
$ws.Range('E33').Value = '  +0.31%  '

# Row 34
$ws.Range('E34').Value = '  +2.41%  '

# Row 35
Set-TextCellValue $ws 'D35' '0.0798'
$ws.Range('E35').Value = '  +4.73%  '

# Row 36
$ws.Range('E36').Value = '  +0.31%  '

# Row 37
$ws.Range('E37').Value = '  +5.95%  '

# Row 38
Set-TextCellValue $ws 'D38' '4.76'
$ws.Range('E38').Value = '  +4.51%  '

# Row 39
Set-TextCellValue $ws 'D39' '3.03'
$ws.Range('E39').Value = '  +3.21%  '

# Row 41
Set-TextCellValue $ws 'D41' '122.98'
$ws.Range('E41').Value = '  -1.49%  '

# Row 42
$ws.Range('E42').Value = '  -0.82%  '

# Row 43
Set-TextCellValue $ws 'D43' '21.38'
$ws.Range('E43').Value = '  +1.52%  '

# Row 44
$ws.Range('E44').Value = '  +2.08%  '

# Row 45
$ws.Range('D45').Value = '2.006.52'
$ws.Range('E45').Value = '  +2.23%  '

# Row 46
Set-TextCellValue $ws 'D46' '3.13'
$ws.Range('E46').Value = '  +5.67%  '

# Row 47
Set-TextCellValue $ws 'D47' '2.08'
$ws.Range('E47').Value = '  -1.69%  '

# Row 48
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCellValue $ws 'D48' '9.12'
$ws.Range('E48').Value = '  +0.48%  '

# Row 49
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCellValue $ws 'D49' '1.78'
$ws.Range('E49').Value = '  -3.48%  '

# Row 50
Set-TextCellValue $ws 'D50' '5.28'
$ws.Range('E50').Value = '  +7.33%  '

# Row 51
Set-TextCellValue $ws 'D51' '78.67'
$ws.Range('E51').Value = '  +1.29%  '

